$p = $ppt.ActivePresentation
try {
    $p.ApplyTheme("ppt/theme/theme1.xml")
    Write-Host "ApplyTheme ok"
} catch {
    Write-Host "ERROR1:" $_.Exception.Message
}

$s = $p.Slides.Item(5)
try {
    $s.ApplyTheme("ppt/theme/theme1.xml")
    Write-Host "Slide.ApplyTheme ok"
} catch {
    Write-Host "ERROR2:" $_.Exception.Message
}
